# Add two new "security" user rows to the Users sheet and tag row 5 with a Role.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 (existing user "abc") gains a Role value in column D.
$ws.Range("A5").Value2 = 4
$ws.Range("B5").Value2 = 'abc'
$ws.Range("C5").Value2 = '$2a$12$4tuF4ECM63Ax3dEqBECu/.HsGl6nvUB7U0qs/nqRNzMv26i.6giWS'
$ws.Range("D5").Value2 = 'ADMIN,STAFF'

# New row 6: user "pqr"
$ws.Range("A6").Value2 = 5
$ws.Range("B6").Value2 = 'pqr'
$ws.Range("C6").Value2 = '$2a$12$Z6brAIjlD7Tu6/3ST2c3aO/M4vpE40UWRIM8wQIfTxq76De7bMhku'

# New row 7: user "mno" with Role "Default"
$ws.Range("A7").Value2 = 5
$ws.Range("B7").Value2 = 'mno'
$ws.Range("C7").Value2 = '$2a$12$VSmfazM8wxdzTyNEnEF6pOD5lAcreSIGTsrkVtnYv382ZBlJRoaQ6'
$ws.Range("D7").Value2 = 'Default'

# Match the author's final cell selection.
$ws.Range("C9").Select()
